# This script reproduces the commit "Added logic for getting information from
# each item and save all the information in Excel Workbook": the worksheet
# previously listed 30 scraped article titles (rows 2-31) under the
# "articles" header in column A. After the scraping/processing logic ran,
# only the first article's info was kept (rows 1-2), the remaining rows
# being removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Group the rows that are about to be removed (this is what leaves the
# outlineLevelRow="1" high-water mark on sheetFormatPr while the
# surviving rows 1-2 stay free of any outlineLevel attribute).
$ws.Rows("3:31").Group()

# Remove all the article rows except the first one, keeping only the
# header ("articles") and the first article ("Mini portable projector").
$ws.Range("A3:A31").EntireRow.Delete()

# Move/restore the active selection like it was left after the edit.
$ws.Range("J10").Select()
